# Applies the "Add leaderboard route and teams data; refactor answer
# validation logic" edit to the Teams Data workbook:
#   1. Removes the "IMPOSTORS" team row (original row 24), shifting the
#      Team Name / Users / Powerups / Score columns of the rows below it
#      up by one.
#   2. Renumbers the Serial No column sequentially.
#   3. Populates the Phase Order / Phase 1-3 Task Order columns (D:G) for
#      every remaining team with a repeating rotation pattern.
#   4. Restores the Credit Card No column (I) to the first 27 values of
#      the original (pre-edit) list - i.e. it does not shift with the
#      row deletion, it simply loses its final entry.
#   5. Applies two standalone Powerups/Score corrections (BABLU, Team
#      Holmes).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 0: remember the original Credit Card No values (I2:I29, 28
# rows) before anything is shifted around. ---
$originalCards = @(
  "MLSC274581924053","MLSC273411206789","MLSC278956012348","MLSC271900439281",
  "MLSC276753908823","MLSC278021677349","MLSC279188325690","MLSC274012093948",
  "MLSC273665718204","MLSC279937456132","MLSC275302947685","MLSC272490411236",
  "MLSC278386074821","MLSC271219486573","MLSC275630089147","MLSC277953712340",
  "MLSC273519849023","MLSC272764021980","MLSC278241857304","MLSC279401358492",
  "MLSC275146789013","MLSC273805276149","MLSC271729503826","MLSC276089314578",
  "MLSC274920348612","MLSC273276041398","MLSC275490028347","MLSC277150283904"
)

# --- Step 1: delete the "IMPOSTORS" team row (row 24). This shifts all
# rows below it up by one for every column. ---
$ws.Rows(24).Delete()

# --- Step 2: renumber Serial No (A2:A28) sequentially. ---
for ($i = 0; $i -lt 27; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $i + 1
}

# --- Step 3: fill Phase Order (D) / Phase 1-3 Task Order (E:G) with the
# rotating pattern, for every data row (2..28). ---
$dPatterns = @("1, 2, 3", "2, 3, 1", "3, 1, 2")
$ePatterns = @(
  "101, 102, 103, 104, 105", "102, 103, 104, 105, 101", "103, 104, 105, 101, 102",
  "104, 105, 101, 102, 103", "105, 101, 102, 103, 104"
)
$fPatterns = @(
  "201, 202, 203, 204, 205", "202, 203, 204, 205, 201", "203, 204, 205, 201, 202",
  "204, 205, 201, 202, 203", "205, 201, 202, 203, 204"
)
$gPatterns = @(
  "301, 302, 303, 304, 305", "302, 303, 304, 305, 301", "303, 304, 305, 301, 302",
  "304, 305, 301, 302, 303", "305, 301, 302, 303, 304"
)

for ($i = 0; $i -lt 27; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 4).Value = $dPatterns[$i % 3]
    $ws.Cells.Item($row, 5).Value = $ePatterns[$i % 5]
    $ws.Cells.Item($row, 6).Value = $fPatterns[$i % 5]
    $ws.Cells.Item($row, 7).Value = $gPatterns[$i % 5]
}

# --- Step 4: restore Credit Card No (I2:I28) to the first 27 of the
# original values (it does not follow the row-24 deletion shift). ---
for ($i = 0; $i -lt 27; $i++) {
    $ws.Cells.Item($i + 2, 9).Value = $originalCards[$i]
}

# --- Step 5: standalone Powerups / Score corrections. ---
# BABLU is now on row 16.
$ws.Cells.Item(16, 8).Value = "5, 2, 7"
$ws.Cells.Item(16, 10).Value = 401

# Team Holmes is now on row 22.
$ws.Cells.Item(22, 8).Value = "2"
$ws.Cells.Item(22, 10).Value = 369
